$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 03:52"

function Set-CountryRow {
    param($row, $values)
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

# Refreshed case numbers caused several countries to re-sort by total cases
# (column B, descending). Write the new country + stats for each affected row.

# Rows 43-46: Mexico / Panama move above Finlandia / Serbia
Set-CountryRow 43 @("Mexico", 2143, 253, 633, 1416, 1, 15, 94)
Set-CountryRow 44 @("Panama", 1988, 187, 13, 1921, 78, 8, 54)
Set-CountryRow 45 @("Finlandia", 1927, 0, 300, 1599, 73, 0, 28)
Set-CountryRow 46 @("Serbia", 1908, 0, 54, 1803, 98, 0, 51)

# Rows 60-61: Nueva Zelanda moves above Estonia
Set-CountryRow 60 @("Nueva Zelanda", 1106, 67, 176, 929, 3, 0, 1)
Set-CountryRow 61 @("Estonia", 1097, 0, 62, 1020, 17, 0, 15)

# Rows 96-97: Honduras moves above Oman
Set-CountryRow 96 @("Honduras", 298, 30, 6, 270, 10, 0, 22)
Set-CountryRow 97 @("Oman", 298, 0, 61, 235, 3, 0, 2)

# Rows 124-126: Paraguay moves above Ruanda / Trinidad yTobago
Set-CountryRow 124 @("Paraguay", 113, 9, 12, 98, 2, 0, 3)
Set-CountryRow 125 @("Ruanda", 104, 0, 4, 100, 0, 0, 0)
Set-CountryRow 126 @("Trinidad yTobago", 104, 0, 1, 96, 0, 0, 7)
